# Append newly completed model-run results (LogisticRegression, ANN, LSTM)
# as rows 39-46 to the results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(39, 1).Value = 'LogisticRegression'
$ws.Cells.Item(39, 2).Value = '{''learning_rate'': 0.001, ''epochs'': 10, ''batch_size'': 32, ''optimizer_type'': ''adam'', ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(39, 3).Value = 0.5397221446037292
$ws.Cells.Item(39, 4).Value = 0.8218031525611877
$ws.Cells.Item(39, 5).Value = 0.5395736694335938
$ws.Cells.Item(39, 6).Value = 0.8220245242118835

$ws.Cells.Item(40, 1).Value = 'ANN'
$ws.Cells.Item(40, 2).Value = '{''hidden_layers'': [32], ''dropout_rate'': 0.3, ''learning_rate'': 0.01, ''epochs'': 20, ''batch_size'': 64, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(40, 3).Value = 0.4079259037971497
$ws.Cells.Item(40, 4).Value = 0.8522424101829529
$ws.Cells.Item(40, 5).Value = 0.3547864258289337
$ws.Cells.Item(40, 6).Value = 0.8685131072998047

$ws.Cells.Item(41, 1).Value = 'ANN'
$ws.Cells.Item(41, 2).Value = '{''hidden_layers'': [32], ''dropout_rate'': 0.3, ''learning_rate'': 0.01, ''epochs'': 20, ''batch_size'': 64, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(41, 3).Value = 0.4055112302303314
$ws.Cells.Item(41, 4).Value = 0.8511757850646973
$ws.Cells.Item(41, 5).Value = 0.3556340336799622
$ws.Cells.Item(41, 6).Value = 0.8694465160369873

$ws.Cells.Item(42, 1).Value = 'LSTM'
$ws.Cells.Item(42, 2).Value = '{''time_steps'': 8, ''lstm_units'': 64, ''epochs'': 10, ''batch_size'': 90, ''learning_rate'': 0.001, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(42, 3).Value = 0.2662752270698547
$ws.Cells.Item(42, 4).Value = 0.9027020931243896
$ws.Cells.Item(42, 5).Value = 0.2909113466739655
$ws.Cells.Item(42, 6).Value = 0.8963072299957275

$ws.Cells.Item(43, 1).Value = 'ANN'
$ws.Cells.Item(43, 2).Value = '{''hidden_layers'': [32], ''dropout_rate'': 0.3, ''learning_rate'': 0.01, ''epochs'': 20, ''batch_size'': 64, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(43, 3).Value = 0.4043055474758148
$ws.Cells.Item(43, 4).Value = 0.8534204959869385
$ws.Cells.Item(43, 5).Value = 0.3485387563705444
$ws.Cells.Item(43, 6).Value = 0.8739444017410278

$ws.Cells.Item(44, 1).Value = 'LSTM'
$ws.Cells.Item(44, 2).Value = '{''time_steps'': 8, ''lstm_units'': 64, ''epochs'': 10, ''batch_size'': 90, ''learning_rate'': 0.001, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(44, 3).Value = 0.2613523900508881
$ws.Cells.Item(44, 4).Value = 0.9052728414535522
$ws.Cells.Item(44, 5).Value = 0.2875173687934875
$ws.Cells.Item(44, 6).Value = 0.8960378170013428

$ws.Cells.Item(45, 1).Value = 'LSTM'
$ws.Cells.Item(45, 2).Value = '{''time_steps'': 8, ''hidden_size'': 128, ''num_layers'': 3, ''dropout_rate'': 0.4, ''epochs'': 10, ''batch_size'': 90, ''learning_rate'': 0.0005, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(45, 3).Value = 0.3039801716804504
$ws.Cells.Item(45, 4).Value = 0.8872016668319702
$ws.Cells.Item(45, 5).Value = 0.3051522970199585
$ws.Cells.Item(45, 6).Value = 0.8897278904914856

$ws.Cells.Item(46, 1).Value = 'LSTM'
$ws.Cells.Item(46, 2).Value = '{''time_steps'': 8, ''lstm_units'': 64, ''epochs'': 10, ''batch_size'': 90, ''learning_rate'': 0.001, ''early_stopping'': True, ''patience'': 10, ''learning_rate_scheduling'': True, ''factor'': 0.1, ''min_lr'': 1e-06}'
$ws.Cells.Item(46, 3).Value = 0.2628472447395325
$ws.Cells.Item(46, 4).Value = 0.9046638011932373
$ws.Cells.Item(46, 5).Value = 0.2815037965774536
$ws.Cells.Item(46, 6).Value = 0.8970574736595154

Write-Host "Appended rows 39-46 to $($ws.Name); new dimension should be A1:F46"
